# DaySale report: two new sale lines were recorded (TELFAST and VOLTAREN)
# between the previous report run and this one, so the item table grows
# from 3 to 5 rows and the trailing "total" / footer rows shift down by 2.
# The footer timestamp is also refreshed to the new generation time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteAll     = -4104
$xlPasteFormats = -4122

function Copy-Row($srcRow, $dstRow) {
    $ws.Range("A" + $srcRow + ":Q" + $srcRow).Copy()
    $ws.Range("A" + $dstRow + ":Q" + $dstRow).PasteSpecial($xlPasteAll)
}

function Copy-RowFormat($srcRow, $dstRow) {
    $ws.Range("A" + $srcRow + ":Q" + $srcRow).Copy()
    $ws.Range("A" + $dstRow + ":Q" + $dstRow).PasteSpecial($xlPasteFormats)
}

# Cells in columns L (qty, custom numeric format) and P (sale price, 0.00
# format) store their numbers as literal text in this report. Typing a
# leading apostrophe forces text entry; re-applying the column's own
# format afterward snaps the cell back onto the shared style used by the
# rest of the table (so we are not left pointing at a throwaway style).
function Set-TextValue($addr, $formatTemplateAddr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($formatTemplateAddr).Copy()
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}

# Plain text columns (C name, H balance, N price, Q ratio) already carry a
# Text number format, so a normal value assignment keeps them as text.
function Set-PlainText($addr, $text) {
    $ws.Range($addr).Value = $text
}

# --- 1. Push the footer row (A11:Q11) down to row 13 and refresh the
#        printed timestamp it carries. ---
Copy-Row 11 13
Set-PlainText "A13" "Monday, 8 September, 2025 10:28 AM"

# --- 2. Push the total row (A10:Q10) down to row 12 and update the sum
#        shown in P12 (old total 84.045 + 160 (TELFAST) + 51 (VOLTAREN)). ---
Copy-Row 10 12
$ws.Range("P12").Value = 295.045

# --- 3. Push the existing item-3 row (TUSSKAN, currently row 9) down to
#        row 10, keeping its values/style untouched. ---
Copy-Row 9 10

# --- 4. Build the new item-5 row (VOLTAREN) in row 11, based on the
#        item-row template style from row 9. ---
Copy-RowFormat 9 11
$ws.Range("A11").Value = 5
Set-PlainText "C11" "VOLTAREN 75MG/3ML 3 AMP."
Set-PlainText "H11" "3:0"
Set-TextValue "L11" "L7" "1"
Set-PlainText "N11" "51.00"
Set-TextValue "P11" "P7" "51.0000"
Set-PlainText "Q11" "1:0"

# --- 5. Turn row 9 into the new item-3 row (TELFAST), overwriting the
#        values that used to belong to TUSSKAN (now safely duplicated in
#        row 10). ---
Set-PlainText "C9" "TELFAST 180MG 20 F.C. TABS"
Set-PlainText "H9" "0:1"
Set-TextValue "L9" "L7" "1"
Set-PlainText "N9" "160.00"
Set-TextValue "P9" "P7" "160.0000"
Set-PlainText "Q9" "1:0"

# --- 6. Re-create the merged cells for the two new item rows (9 already
#        has its original merges; 10/11 need the same A:B / C:G / H:K /
#        L:M / N:O layout used by every item row). ---
$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

$ws.Range("A11:B11").Merge()
$ws.Range("C11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()
$ws.Range("N11:O11").Merge()

Write-Output "DaySale report updated with TELFAST and VOLTAREN sale rows"
